$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text updates (Volume/Number and the reporting week date range)
#    A8  = "Volume 31   Number  35" -> "...36"
#    C9  = "Report Covering the Week  8/26/2024  Through  9/1/2024"
#           -> "...9/2/2024  Through  9/8/2024"
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/2/2024  Through  9/8/2024"

# ---------------------------------------------------------------------------
# 2) Cells whose content switches between a numeric value and the special
#    text markers ("0" / "***.*") used for N/A counts and percentages.
#    Copying a same-shaped cell (value + number format) from elsewhere on
#    the sheet reproduces both the literal text/number and the exact style.
# ---------------------------------------------------------------------------
$typeSwitches = [ordered]@{
    "D14" = "C14"   # -> text "0"
    "E14" = "L14"   # -> text "***.*"
    "C18" = "F14"   # -> text "0"
    "D18" = "F22"   # -> number 2
    "E18" = "H14"   # -> number -100
    "C20" = "I29"   # -> number 2
    "C22" = "D22"   # -> text "0"
    "C23" = "J29"   # -> number 2
    "D23" = "C15"   # -> text "0"
    "E23" = "E22"   # -> text "***.*"
    "D28" = "D15"   # -> number 1
    "E28" = "H27"   # -> number -100
    "F28" = "C28"   # -> text "0"
    "D31" = "C29"   # -> text "0"
    "E31" = "H22"   # -> text "***.*"
}
foreach ($dest in $typeSwitches.Keys) {
    $src = $typeSwitches[$dest]
    $ws.Range($src).Copy($ws.Range($dest))
}

# ---------------------------------------------------------------------------
# 3) Plain numeric refreshes (counts + computed % changes) across rows 15-31
# ---------------------------------------------------------------------------
$valueUpdates = @{
    "G15" = 2
    "J15" = 8
    "K15" = -12.5
    "C16" = 1
    "D16" = 2
    "E16" = -50
    "F16" = 8
    "G16" = 4
    "H16" = 100
    "I16" = 69
    "J16" = 62
    "K16" = 11.290322580645
    "L16" = -14.814814814814
    "M16" = 4.545454545454
    "N16" = -83.211678832116
    "C17" = 2
    "D17" = 1
    "E17" = 100
    "F17" = 9
    "G17" = 12
    "H17" = -25
    "I17" = 82
    "J17" = 80
    "K17" = 2.5
    "L17" = 15.492957746478
    "M17" = 78.260869565217
    "N17" = 24.242424242424
    "F18" = 3
    "G18" = 7
    "H18" = -57.142857142857
    "I18" = 92
    "J18" = 81
    "K18" = 13.580246913580
    "L18" = -14.018691588785
    "M18" = 22.666666666666
    "N18" = -84.893267651888
    "C19" = 20
    "D19" = 16
    "E19" = 25
    "F19" = 57
    "G19" = 56
    "H19" = 1.785714285714
    "I19" = 529
    "J19" = 502
    "K19" = 5.378486055776
    "L19" = -6.205673758865
    "M19" = 18.080357142857
    "N19" = -60.374531835206
    "D20" = 3
    "E20" = -33.333333333333
    "G20" = 4
    "H20" = 125
    "I20" = 36
    "J20" = 73
    "K20" = -50.684931506849
    "M20" = 71.428571428571
    "N20" = -95.342820181112
    "C21" = 25
    "D21" = 25
    "E21" = 0
    "F21" = 86
    "G21" = 88
    "H21" = -2.272727272727
    "I21" = 815
    "J21" = 810
    "K21" = 0.617283950617
    "L21" = -8.220720720720
    "M21" = 22.740963855421
    "N21" = -74.634298163709
    "F23" = 4
    "H23" = -20
    "I23" = 29
    "K23" = -6.451612903225
    "L23" = 52.631578947368
    "M23" = 38.095238095238
    "C24" = 21
    "D24" = 33
    "E24" = -36.363636363636
    "F24" = 93
    "G24" = 124
    "H24" = -25
    "I24" = 879
    "J24" = 801
    "K24" = 9.737827715355
    "L24" = -4.560260586319
    "M24" = 19.267299864314
    "C25" = 21
    "D25" = 20
    "E25" = 5
    "F25" = 69
    "G25" = 78
    "H25" = -11.538461538461
    "I25" = 691
    "J25" = 604
    "K25" = 14.403973509933
    "L25" = -6.621621621621
    "C26" = 6
    "D26" = 6
    "E26" = 0
    "F26" = 12
    "G26" = 17
    "H26" = -29.411764705882
    "I26" = 163
    "J26" = 159
    "K26" = 2.515723270440
    "L26" = 12.413793103448
    "M26" = -15.104166666666
    "G27" = 2
    "J27" = 16
    "K27" = -43.75
    "H28" = -100
    "J28" = 28
    "K28" = -35.714285714285
    "L28" = -58.139534883720
    "G31" = 1
}
foreach ($ref in $valueUpdates.Keys) {
    $ws.Range($ref).Value = $valueUpdates[$ref]
}
